$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Metadata") - update version/status/experimental/date/description
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item(1)

# Version: 1.0.0 -> 0.1.0
$meta.Range("B3").Value = "0.1.0"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Experimental: (blank) -> false
# A bare "false" is auto-typed as a Boolean cell by the engine, but the
# target file stores it as literal text ("t=s"). Use a leading apostrophe
# to force text entry, then restore the surrounding body style (format-only
# paste from a same-styled cell) so the quote-prefix bookkeeping the
# apostrophe triggers doesn't linger on the cell itself.
$meta.Range("B7").Value = "'false"
$meta.Range("B9").Copy()
$meta.Range("B7").PasteSpecial(-4122)

# Date: refreshed timestamp
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: (blank) -> new description text
$meta.Range("B11").Value = "Value set for evaluating patient goal outcomes"

# ---------------------------------------------------------------------------
# Sheet 2 ("Include #0") - update/add SNOMED concept rows
# ---------------------------------------------------------------------------
$inc = $wb.Worksheets.Item(2)

# Row 2: concept code changes from 385633008 to 370996005 (kept as text,
# same apostrophe + format-restore trick as above since it is all digits
# and would otherwise be stored as a number), and gets a description
# ("Patient condition resolved") that it didn't have before.
$inc.Range("A2").Value = "'370996005"
$inc.Range("A4").Copy()
$inc.Range("A2").PasteSpecial(-4122)
$inc.Range("B2").Value = "Patient condition resolved"

# Row 3: concept code 385634002 is unchanged (leave it alone), but now gets
# a description ("Worsened") that it didn't have before.
$inc.Range("B3").Value = "Worsened"

# Insert a brand-new row at position 4 (pushing the old blank row 4 and the
# "System URI" row down by one), then restore the formatting from the row
# above it so the new cells pick up the same border/alignment style used by
# the rest of the table body.
$inc.Rows.Item(4).Insert()
$inc.Range("A3:B3").Copy()
$inc.Range("A4:B4").PasteSpecial(-4122)

# New row 4: concept 118222006 / "General finding of observation of patient"
$inc.Range("A4").Value = "'118222006"
$inc.Range("A3").Copy()
$inc.Range("A4").PasteSpecial(-4122)
$inc.Range("B4").Value = "General finding of observation of patient"
